$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1715328467153285
$ws.Range("C2").Value = 0.5948905109489051
$ws.Range("J2").Value = 0.01094890510948905
$ws.Range("P2").Value = 0.1167883211678832
$ws.Range("S2").Value = 0.1058394160583942
$ws.Range("C3").Value = 0.04191616766467066
$ws.Range("J3").Value = 0.02994011976047904
$ws.Range("P3").Value = 0.7065868263473054
$ws.Range("S3").Value = 0.2215568862275449
$ws.Range("J4").Value = 0.05882352941176471
$ws.Range("P4").Value = 0.7058823529411765
$ws.Range("S4").Value = 0.2352941176470588
$ws.Range("P5").Value = 1
$ws.Range("B6").Value = 0.06572769953051644
$ws.Range("F6").Value = 0.02816901408450704
$ws.Range("J6").Value = 0.2394366197183098
$ws.Range("O6").Value = 0.0187793427230047
$ws.Range("Q6").Value = 0.2394366197183098
$ws.Range("R6").Value = 0.08450704225352113
$ws.Range("S6").Value = 0.323943661971831
$ws.Range("B7").Value = 0.09941520467836257
$ws.Range("D7").Value = 0.005847953216374269
$ws.Range("F7").Value = 0.02923976608187134
$ws.Range("J7").Value = 0.09941520467836257
$ws.Range("O7").Value = 0.01169590643274854
$ws.Range("Q7").Value = 0.2046783625730994
$ws.Range("R7").Value = 0.08771929824561403
$ws.Range("S7").Value = 0.4619883040935672
$ws.Range("B8").Value = 0.0765661252900232
$ws.Range("D8").Value = 0.009280742459396751
$ws.Range("F8").Value = 0.06264501160092807
$ws.Range("J8").Value = 0.1183294663573086
$ws.Range("O8").Value = 0.01624129930394431
$ws.Range("Q8").Value = 0.2320185614849188
$ws.Range("R8").Value = 0.08816705336426914
$ws.Range("S8").Value = 0.3967517401392112
$ws.Range("B9").Value = 0.111587982832618
$ws.Range("D9").Value = 0.008583690987124463
$ws.Range("F9").Value = 0.1030042918454936
$ws.Range("J9").Value = 0.1030042918454936
$ws.Range("O9").Value = 0.0128755364806867
$ws.Range("Q9").Value = 0.2188841201716738
$ws.Range("R9").Value = 0.09871244635193133
$ws.Range("S9").Value = 0.3433476394849785
$ws.Range("B10").Value = 0.1127241673783091
$ws.Range("D10").Value = 0.02305721605465414
$ws.Range("E10").Value = 0.001707941929974381
$ws.Range("F10").Value = 0.06917164816396243
$ws.Range("J10").Value = 0.08625106746370624
$ws.Range("O10").Value = 0.01878736122971819
$ws.Range("Q10").Value = 0.2356959863364645
$ws.Range("R10").Value = 0.08710503842869342
$ws.Range("S10").Value = 0.3654995730145175
$ws.Range("G11").Value = 0.09126984126984126
$ws.Range("J11").Value = 0.09126984126984126
$ws.Range("K11").Value = 0.1587301587301587
$ws.Range("L11").Value = 0.623015873015873
$ws.Range("S11").Value = 0.03571428571428571
$ws.Range("G12").Value = 0.8024691358024691
$ws.Range("J12").Value = 0.1234567901234568
$ws.Range("K12").Value = 0.01234567901234568
$ws.Range("L12").Value = 0.02469135802469136
$ws.Range("S12").Value = 0.03703703703703703
$ws.Range("G13").Value = 0.6285714285714286
$ws.Range("J13").Value = 0.3714285714285714
$ws.Range("G14").Value = 0.6666666666666666
$ws.Range("S14").Value = 0.3333333333333333
$ws.Range("F15").Value = 0.01923076923076923
$ws.Range("H15").Value = 0.1634615384615385
$ws.Range("I15").Value = 0.0625
$ws.Range("J15").Value = 0.3605769230769231
$ws.Range("K15").Value = 0.0625
$ws.Range("M15").Value = 0.01442307692307692
$ws.Range("N15").Value = 0.004807692307692308
$ws.Range("O15").Value = 0.07211538461538461
$ws.Range("S15").Value = 0.2403846153846154
$ws.Range("F16").Value = 0.04191616766467066
$ws.Range("H16").Value = 0.2215568862275449
$ws.Range("I16").Value = 0.1077844311377246
$ws.Range("J16").Value = 0.3473053892215569
$ws.Range("K16").Value = 0.1137724550898204
$ws.Range("M16").Value = 0.01796407185628742
$ws.Range("O16").Value = 0.02395209580838323
$ws.Range("S16").Value = 0.125748502994012
$ws.Range("F17").Value = 0.01330798479087452
$ws.Range("H17").Value = 0.1615969581749049
$ws.Range("I17").Value = 0.1083650190114068
$ws.Range("J17").Value = 0.4106463878326996
$ws.Range("K17").Value = 0.08555133079847908
$ws.Range("M17").Value = 0.01520912547528517
$ws.Range("N17").Value = 0.001901140684410646
$ws.Range("O17").Value = 0.06844106463878327
$ws.Range("S17").Value = 0.1349809885931559
$ws.Range("F18").Value = 0.01015228426395939
$ws.Range("H18").Value = 0.1878172588832487
$ws.Range("I18").Value = 0.1472081218274112
$ws.Range("J18").Value = 0.4467005076142132
$ws.Range("K18").Value = 0.08629441624365482
$ws.Range("M18").Value = 0.01522842639593909
$ws.Range("O18").Value = 0.05583756345177665
$ws.Range("S18").Value = 0.05076142131979695
$ws.Range("F19").Value = 0.0247229326513214
$ws.Range("H19").Value = 0.2020460358056266
$ws.Range("I19").Value = 0.0988917306052856
$ws.Range("J19").Value = 0.3742540494458653
$ws.Range("K19").Value = 0.0937766410912191
$ws.Range("M19").Value = 0.01705029838022165
$ws.Range("N19").Value = 0.0008525149190110827
$ws.Range("O19").Value = 0.06820119352088662
$ws.Range("S19").Value = 0.1202046035805627
